$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1679.747
$ws.Range("I40").Value = 1665.7377
$ws.Range("J40").Value = 1718.591
$ws.Range("K40").Value = 1665.7377
$ws.Range("L40").Value = 1718.591
$ws.Range("M40").Value = -1490.7377
$ws.Range("N40").Value = -2068.591

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3728.0881
$ws.Range("I64").Value = 3729.3794
$ws.Range("J64").Value = 3720.6
$ws.Range("K64").Value = 3729.3794
$ws.Range("L64").Value = 3720.6
$ws.Range("M64").Value = -3481.3794
$ws.Range("N64").Value = -4216.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 3728.0881
$ws.Range("I67").Value = 3729.3794
$ws.Range("J67").Value = 3720.6
$ws.Range("K67").Value = 3729.3794
$ws.Range("L67").Value = 3720.6
$ws.Range("M67").Value = -2871.3794
$ws.Range("N67").Value = -5436.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1384.16
$ws.Range("I132").Value = 1384.16
$ws.Range("K132").Value = 4152.48
$ws.Range("M132").Value = -1622.48

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1396.2413
$ws.Range("I137").Value = 1348.0435
$ws.Range("K137").Value = 4044.1305
$ws.Range("M137").Value = -1494.1305

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5061.8525
$ws.Range("I32").Value = 3761.0833
$ws.Range("K32").Value = 3761.0833
$ws.Range("M32").Value = -3474.0833

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 8189.5293
$ws.Range("I61").Value = 8513.875
$ws.Range("K61").Value = 8513.875
$ws.Range("M61").Value = -8301.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1452.5857
$ws.Range("I74").Value = 1438.625
$ws.Range("K74").Value = 1438.625
$ws.Range("M74").Value = -564.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1452.5857
$ws.Range("I77").Value = 1438.625
$ws.Range("K77").Value = 7193.125
$ws.Range("M77").Value = -2825.125

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2367.1667
$ws.Range("I88").Value = 2299
$ws.Range("K88").Value = 2299
$ws.Range("M88").Value = -1893

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H91").Value = 2367.1667
$ws.Range("I91").Value = 2299
$ws.Range("K91").Value = 2299
$ws.Range("M91").Value = -895

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2968.7317
$ws.Range("I132").Value = 1303.7778
$ws.Range("J132").Value = 6179.7144
$ws.Range("K132").Value = 3911.3334
$ws.Range("L132").Value = 18539.1432
$ws.Range("M132").Value = -1381.3334
$ws.Range("N132").Value = -23599.1432

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 8189.5293
$ws.Range("I136").Value = 8513.875
$ws.Range("K136").Value = 25541.625
$ws.Range("M136").Value = -22991.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 516.125
$ws.Range("I22").Value = 525
$ws.Range("J22").Value = 507.25
$ws.Range("K22").Value = 525
$ws.Range("L22").Value = 507.25
$ws.Range("M22").Value = -352
$ws.Range("N22").Value = -853.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 41668308
$ws.Range("I86").Value = 47620492
$ws.Range("J86").Value = 3000
$ws.Range("K86").Value = 47620492
$ws.Range("L86").Value = 3000
$ws.Range("M86").Value = -47619369
$ws.Range("N86").Value = -5246

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 41668308
$ws.Range("I89").Value = 47620492
$ws.Range("J89").Value = 3000
$ws.Range("K89").Value = 238102460
$ws.Range("L89").Value = 15000
$ws.Range("M89").Value = -238096844
$ws.Range("N89").Value = -26232

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5841.3
$ws.Range("I134").Value = 7247.2
$ws.Range("K134").Value = 21741.6
$ws.Range("M134").Value = -19206.6

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5151.878
$ws.Range("I31").Value = 1789
$ws.Range("J31").Value = 9045.736999999999
$ws.Range("K31").Value = 1789
$ws.Range("L31").Value = 9045.736999999999
$ws.Range("M31").Value = -1494
$ws.Range("N31").Value = -9635.736999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 5151.878
$ws.Range("I34").Value = 1789
$ws.Range("J34").Value = 9045.736999999999
$ws.Range("K34").Value = 1789
$ws.Range("L34").Value = 9045.736999999999
$ws.Range("M34").Value = -1587
$ws.Range("N34").Value = -9449.736999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 5521.8667
$ws.Range("I62").Value = 6406.7144
$ws.Range("J62").Value = 4747.625
$ws.Range("K62").Value = 6406.7144
$ws.Range("L62").Value = 4747.625
$ws.Range("M62").Value = -5782.7144
$ws.Range("N62").Value = -5995.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 5521.8667
$ws.Range("I65").Value = 6406.7144
$ws.Range("J65").Value = 4747.625
$ws.Range("K65").Value = 32033.572
$ws.Range("L65").Value = 23738.125
$ws.Range("M65").Value = -28913.572
$ws.Range("N65").Value = -29978.125

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H87").Value = 22000
$ws.Range("J87").Value = 22000
$ws.Range("L87").Value = 22000
$ws.Range("N87").Value = -24372

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H90").Value = 22000
$ws.Range("J90").Value = 22000
$ws.Range("L90").Value = 66000
$ws.Range("N90").Value = -77856

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2766.5588
$ws.Range("I134").Value = 3277.1738
$ws.Range("J134").Value = 1698.909
$ws.Range("K134").Value = 9831.5214
$ws.Range("L134").Value = 5096.727000000001
$ws.Range("M134").Value = -7296.5214
$ws.Range("N134").Value = -10166.727

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 316646.75
$ws.Range("I5").Value = 508
$ws.Range("J5").Value = 751337.5
$ws.Range("K5").Value = 1524
$ws.Range("L5").Value = 2254012.5
$ws.Range("M5").Value = -1412
$ws.Range("N5").Value = -2254236.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H127").Value = 15151959
$ws.Range("J127").Value = 15151959
$ws.Range("L127").Value = 45455877
$ws.Range("N127").Value = -45465797

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 10553.617
$ws.Range("I134").Value = 12891.111
$ws.Range("K134").Value = 38673.333
$ws.Range("M134").Value = -33603.333

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H135").Value = 316646.75
$ws.Range("I135").Value = 508
$ws.Range("J135").Value = 751337.5
$ws.Range("K135").Value = 4572
$ws.Range("L135").Value = 6762037.5
$ws.Range("M135").Value = -2037
$ws.Range("N135").Value = -6767107.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3276.6597
$ws.Range("I132").Value = 3035.2307
$ws.Range("J132").Value = 3575.5715
$ws.Range("K132").Value = 9105.6921
$ws.Range("L132").Value = 10726.7145
$ws.Range("M132").Value = -6575.6921
$ws.Range("N132").Value = -15786.7145

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 90911290
$ws.Range("J68").Value = 200002780
$ws.Range("L68").Value = 200002780
$ws.Range("N68").Value = -200004278

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H71").Value = 90911290
$ws.Range("J71").Value = 200002780
$ws.Range("L71").Value = 1000013900
$ws.Range("N71").Value = -1000021388

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 4525913.5
$ws.Range("I122").Value = 4764196.5
$ws.Range("J122").Value = 3334500
$ws.Range("K122").Value = 14292589.5
$ws.Range("L122").Value = 10003500
$ws.Range("M122").Value = -14290139.5
$ws.Range("N122").Value = -10008400

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 27093278
$ws.Range("I132").Value = 39406030
$ws.Range("J132").Value = 5219.6
$ws.Range("K132").Value = 118218090
$ws.Range("L132").Value = 15658.8
$ws.Range("M132").Value = -118215560
$ws.Range("N132").Value = -20718.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 5446.8125
$ws.Range("I136").Value = 5570.9697
$ws.Range("J136").Value = 5173.6665
$ws.Range("K136").Value = 16712.9091
$ws.Range("L136").Value = 15520.9995
$ws.Range("M136").Value = -14162.9091
$ws.Range("N136").Value = -20620.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H57").Value = 53000
$ws.Range("I57").Value = 53000
$ws.Range("K57").Value = 53000
$ws.Range("M57").Value = -52246

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1881
$ws.Range("I132").Value = 1259.5834
$ws.Range("K132").Value = 3778.7502
$ws.Range("M132").Value = -1248.7502

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1461.8462
$ws.Range("I136").Value = 899.8570999999999
$ws.Range("J136").Value = 2117.5
$ws.Range("K136").Value = 2699.5713
$ws.Range("L136").Value = 6352.5
$ws.Range("M136").Value = -149.5712999999996
$ws.Range("N136").Value = -11452.5
